{"js": "// The whole body content is being replaced: the three \"NormalWeb\"-styled\n// Q&A paragraphs (plus the trailing blank paragraph) are swapped out for six\n// brand-new, plain (unstyled) paragraphs:\n//   \"codingmadeclear.com\", \"\", \"mock\", \"lowercase/end in s\", \"\", \"the schema\"\nconst body = context.document.body;\n\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Delete every existing paragraph except the last one. The last paragraph in\n// this document is the pre-existing blank paragraph mark with no direct\n// formatting, so keeping it (instead of clearing+recreating everything) means\n// the new content we type into it picks up no leftover style/run formatting.\nfor (let i = 0; i < paragraphs.items.length - 1; i++) {\n  paragraphs.items[i].delete();\n}\nawait context.sync();\n\n// Reuse that last, formatting-free paragraph for the new first line of text.\nconst firstParagraph = body.paragraphs.getLast();\nfirstParagraph.insertText(\"codingmadeclear.com\", Word.InsertLocation.replace);\n\n// Append the remaining five paragraphs (including the two blank ones) in\n// order at the end of the body.\nbody.insertParagraph(\"\", Word.InsertLocation.end);\nbody.insertParagraph(\"mock\", Word.InsertLocation.end);\nbody.insertParagraph(\"lowercase/end in s\", Word.InsertLocation.end);\nbody.insertParagraph(\"\", Word.InsertLocation.end);\nbody.insertParagraph(\"the schema\", Word.InsertLocation.end);\n\nawait context.sync();\n", "ps1": "# The whole body content is being replaced: the three \"NormalWeb\"-styled\n# Q&A paragraphs (plus the trailing blank paragraph) are swapped out for six\n# brand-new, plain (unstyled) paragraphs:\n#   \"codingmadeclear.com\", \"\", \"mock\", \"lowercase/end in s\", \"\", \"the schema\"\n$d = $word.ActiveDocument\n\n# Delete every existing paragraph except the last one. The last paragraph in\n# this document is the pre-existing blank paragraph mark with no direct\n# formatting, so keeping it (instead of nuking + rebuilding everything) means\n# the new content we type into it picks up no leftover style/run formatting.\nwhile ($d.Paragraphs.Count -gt 1) {\n    $d.Paragraphs(1).Range.Delete()\n}\n\n# Reuse that last, formatting-free paragraph for the new first line of text.\n$d.Paragraphs(1).Range.Text = \"codingmadeclear.com\"\n\n# Append the remaining five paragraphs (including the two blank ones) in\n# order at the end of the document.\n$d.Content.InsertParagraphAfter()\n$d.Content.InsertParagraphAfter()\n$d.Paragraphs($d.Paragraphs.Count).Range.Text = \"mock\"\n\n$d.Content.InsertParagraphAfter()\n$d.Paragraphs($d.Paragraphs.Count).Range.Text = \"lowercase/end in s\"\n\n$d.Content.InsertParagraphAfter()\n$d.Content.InsertParagraphAfter()\n$d.Paragraphs($d.Paragraphs.Count).Range.Text = \"the schema\"\n"}
